$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-removed rows 10-13 (cluster "Resolving-Mac" sender block)
$ws.Rows("10:13").Delete()

# Row 2
$ws.Range("A2").Value = 'FAPs'
$ws.Range("B2").Value = 'Has2'
$ws.Range("C2").Value = 'Hmmr'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 23.544642
$ws.Range("H2").Value = 70.633926
$ws.Range("I2").Value = 0.8903342714957673
$ws.Range("J2").Value = 0.8903342714957674
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.72629
$ws.Range("N2").Value = 17.17887
$ws.Range("O2").Value = 0.3730133452890743
$ws.Range("P2").Value = 0.3730133452890744
$ws.Range("Q2").Value = 134.82344803818
$ws.Range("R2").Value = 1213.41103234362
$ws.Range("S2").Value = 0.3321065650361471
$ws.Range("T2").Value = 0.3321065650361472

# Row 3
$ws.Range("A3").Value = 'FAPs'
$ws.Range("B3").Value = 'Has2'
$ws.Range("C3").Value = 'Hmmr'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 23.544642
$ws.Range("H3").Value = 70.633926
$ws.Range("I3").Value = 0.8903342714957673
$ws.Range("J3").Value = 0.8903342714957674
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.270604
$ws.Range("N3").Value = 3.811812
$ws.Range("O3").Value = 0.08276776911013571
$ws.Range("P3").Value = 0.08276776911013571
$ws.Range("Q3").Value = 29.915916303768
$ws.Range("R3").Value = 269.243246733912
$ws.Range("S3").Value = 0.07369098141400254
$ws.Range("T3").Value = 0.07369098141400256

# Row 4
$ws.Range("A4").Value = 'FAPs'
$ws.Range("B4").Value = 'Has2'
$ws.Range("C4").Value = 'Hmmr'
$ws.Range("D4").Value = 'MuSCs'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 23.544642
$ws.Range("H4").Value = 70.633926
$ws.Range("I4").Value = 0.8903342714957673
$ws.Range("J4").Value = 0.8903342714957674
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.20253
$ws.Range("N4").Value = 21.60759
$ws.Range("O4").Value = 0.4691763445171162
$ws.Range("P4").Value = 0.4691763445171162
$ws.Range("Q4").Value = 169.58099034426
$ws.Range("R4").Value = 1526.22891309834
$ws.Range("S4").Value = 0.4177237788986937
$ws.Range("T4").Value = 0.4177237788986938

# Row 5
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Has2'
$ws.Range("C5").Value = 'Hmmr'
$ws.Range("D5").Value = 'Resolving-Mac'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.544642
$ws.Range("H5").Value = 70.633926
$ws.Range("I5").Value = 0.8903342714957673
$ws.Range("J5").Value = 0.8903342714957674
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.152010666666667
$ws.Range("N5").Value = 3.456032
$ws.Range("O5").Value = 0.07504254108367374
$ws.Range("P5").Value = 0.07504254108367375
$ws.Range("Q5").Value = 27.123678726848
$ws.Range("R5").Value = 244.113108541632
$ws.Range("S5").Value = 0.06681294614692385
$ws.Range("T5").Value = 0.06681294614692386

# Row 6
$ws.Range("A6").Value = 'MuSCs'
$ws.Range("B6").Value = 'Has2'
$ws.Range("C6").Value = 'Hmmr'
$ws.Range("D6").Value = 'ECs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.900079666666667
$ws.Range("H6").Value = 8.700239
$ws.Range("I6").Value = 0.1096657285042327
$ws.Range("J6").Value = 0.1096657285042327
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.72629
$ws.Range("N6").Value = 17.17887
$ws.Range("O6").Value = 0.3730133452890743
$ws.Range("P6").Value = 0.3730133452890744
$ws.Range("Q6").Value = 16.60669719443667
$ws.Range("R6").Value = 149.46027474993
$ws.Range("S6").Value = 0.04090678025292723
$ws.Range("T6").Value = 0.04090678025292724

# Row 7
$ws.Range("A7").Value = 'MuSCs'
$ws.Range("B7").Value = 'Has2'
$ws.Range("C7").Value = 'Hmmr'
$ws.Range("D7").Value = 'FAPs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.900079666666667
$ws.Range("H7").Value = 8.700239
$ws.Range("I7").Value = 0.1096657285042327
$ws.Range("J7").Value = 0.1096657285042327
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.270604
$ws.Range("N7").Value = 3.811812
$ws.Range("O7").Value = 0.08276776911013571
$ws.Range("P7").Value = 0.08276776911013571
$ws.Range("Q7").Value = 3.684852824785333
$ws.Range("R7").Value = 33.163675423068
$ws.Range("S7").Value = 0.009076787696133161
$ws.Range("T7").Value = 0.009076787696133161

# Row 8
$ws.Range("A8").Value = 'MuSCs'
$ws.Range("B8").Value = 'Has2'
$ws.Range("C8").Value = 'Hmmr'
$ws.Range("D8").Value = 'MuSCs'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.900079666666667
$ws.Range("H8").Value = 8.700239
$ws.Range("I8").Value = 0.1096657285042327
$ws.Range("J8").Value = 0.1096657285042327
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.20253
$ws.Range("N8").Value = 21.60759
$ws.Range("O8").Value = 0.4691763445171162
$ws.Range("P8").Value = 0.4691763445171162
$ws.Range("Q8").Value = 20.88791080155667
$ws.Range("R8").Value = 187.99119721401
$ws.Range("S8").Value = 0.05145256561842241
$ws.Range("T8").Value = 0.05145256561842241

# Row 9
$ws.Range("A9").Value = 'MuSCs'
$ws.Range("B9").Value = 'Has2'
$ws.Range("C9").Value = 'Hmmr'
$ws.Range("D9").Value = 'Resolving-Mac'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.900079666666667
$ws.Range("H9").Value = 8.700239
$ws.Range("I9").Value = 0.1096657285042327
$ws.Range("J9").Value = 0.1096657285042327
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.152010666666667
$ws.Range("N9").Value = 3.456032
$ws.Range("O9").Value = 0.07504254108367374
$ws.Range("P9").Value = 0.07504254108367375
$ws.Range("Q9").Value = 3.340922710183111
$ws.Range("R9").Value = 30.068304391648
$ws.Range("S9").Value = 0.008229594936749893
$ws.Range("T9").Value = 0.008229594936749894

